# Update "想去人数" (want-to-go count) figures that changed between scrapes.
# Each worksheet lists events; "全部类型" aggregates the same events found in
# the other three sheets, so the same underlying counts are bumped there too.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1571
$ws1.Range("F5").Value  = 229
$ws1.Range("F7").Value  = 1128
$ws1.Range("F8").Value  = 729
$ws1.Range("F9").Value  = 776
$ws1.Range("F10").Value = 1400
$ws1.Range("F14").Value = 64
$ws1.Range("F17").Value = 441
$ws1.Range("F20").Value = 293

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 995
$ws2.Range("F5").Value = 261
$ws2.Range("F7").Value = 141

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 216

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 216
$ws4.Range("F4").Value  = 1571
$ws4.Range("F7").Value  = 229
$ws4.Range("F8").Value  = 995
$ws4.Range("F10").Value = 1128
$ws4.Range("F11").Value = 729
$ws4.Range("F12").Value = 776
$ws4.Range("F13").Value = 1400
$ws4.Range("F17").Value = 64
$ws4.Range("F20").Value = 441
$ws4.Range("F23").Value = 261
$ws4.Range("F25").Value = 293
$ws4.Range("F27").Value = 141
$ws4.Range("F28").Value = 141
